$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# --- Crime data table updates (weekly section rows 14-30, historical rows 36-42) ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("G14").NumberFormat = "General"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("H14").NumberFormat = "General"
$ws.Range("N14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N14").Value = -100
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C15").NumberFormat = "General"
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = -66.666666666666
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 11
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 6
$ws.Range("L16").Value = -20
$ws.Range("N16").Value = -90.243902439024
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 43.75
$ws.Range("I17").Value = 15
$ws.Range("J17").Value = 7
$ws.Range("K17").Value = 114.285714285714
$ws.Range("L17").Value = 87.5
$ws.Range("M17").Value = 400
$ws.Range("N17").Value = -37.5
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -26.315789473684
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = -12.5
$ws.Range("L18").Value = 40
$ws.Range("M18").Value = -36.363636363636
$ws.Range("N18").Value = -84.444444444444
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -29.411764705882
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -11.111111111111
$ws.Range("I19").Value = 23
$ws.Range("J19").Value = 27
$ws.Range("K19").Value = -14.814814814814
$ws.Range("L19").Value = -54
$ws.Range("M19").Value = 155.555555555556
$ws.Range("N19").Value = 76.923076923076
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 77.777777777777
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 4
$ws.Range("L20").Value = 166.666666666667
$ws.Range("M20").Value = 14.285714285714
$ws.Range("N20").Value = -85.714285714285
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -3.703703703703
$ws.Range("F21").Value = 105
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = 1.941747572815
$ws.Range("I21").Value = 58
$ws.Range("J21").Value = 55
$ws.Range("K21").Value = 5.454545454545
$ws.Range("L21").Value = -18.309859154929
$ws.Range("M21").Value = 56.756756756756
$ws.Range("N21").Value = -67.955801104972
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -50
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -32.258064516129
$ws.Range("I24").Value = 42
$ws.Range("J24").Value = 54
$ws.Range("K24").Value = -22.222222222222
$ws.Range("L24").Value = 10.526315789473
$ws.Range("M24").Value = -4.545454545454
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 57.142857142857
$ws.Range("F25").Value = 38
$ws.Range("H25").Value = 5.555555555555
$ws.Range("I25").Value = 17
$ws.Range("J25").Value = 16
$ws.Range("K25").Value = 6.25
$ws.Range("L25").Value = 88.888888888888
$ws.Range("M25").Value = 0
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C26").NumberFormat = "General"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = -66.666666666666
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C27").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("J27").Value = 1
$ws.Range("K27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K27").Value = 100
$ws.Range("L27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L27").Value = 100
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E28").NumberFormat = "General"
$ws.Range("N28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N28").Value = -50
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C29").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E29").NumberFormat = "General"
$ws.Range("N29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N29").Value = -50
